$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.715.77"
$ws.Range("E2").Value = "  +3.10%  "

$ws.Range("D3").Value = "1.689.74"
$ws.Range("E3").Value = "  +3.51%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5362"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2694"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06449"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07803"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.52%  "

$ws.Range("D12").Value = "1.692.37"
$ws.Range("E12").Value = "  +2.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.522"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5660"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.86%  "

$ws.Range("D15").Value = "0.0₅8508"
$ws.Range("E15").Value = "  +8.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.42%  "

$ws.Range("D17").Value = "26.750.01"
$ws.Range("E17").Value = "  +3.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.830"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "196.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.407"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.03"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1282"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.495"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.415"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06200"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.17%  "

$ws.Range("E30").Value = "  +3.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.616"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.480"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.713"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.020"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.805"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.413"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5741"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.91%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01655"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.976"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.58%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8682"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.54%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.072.36"
$ws.Range("E41").Value = "  +3.95%  "

$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("D44").Value = "1.842.17"
$ws.Range("E44").Value = "  +3.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.21%  "

$ws.Range("E46").Value = "  -1.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.206"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.72%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.107"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4248"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.41%  "
